# Apply scheduled-runner updates to Sheets workbook (Kujata_Profits)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5475
$ws.Range("J64").Value = 3960
$ws.Range("L64").Value = 3960
$ws.Range("N64").Value = -4456
$ws.Range("H67").Value = 5475
$ws.Range("J67").Value = 3960
$ws.Range("L67").Value = 3960
$ws.Range("N67").Value = -5676
$ws.Range("H69").Value = 3349.4443
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 3393.125
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 10179.375
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -11927.375
$ws.Range("H72").Value = 3349.4443
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 3393.125
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 30538.125
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -39274.125
$ws.Range("H98").Value = 9656.053
$ws.Range("I98").Value = 11104.0625
$ws.Range("J98").Value = 1933.3334
$ws.Range("K98").Value = 11104.0625
$ws.Range("L98").Value = 1933.3334
$ws.Range("M98").Value = -9606.0625
$ws.Range("N98").Value = -4929.3334
$ws.Range("H122").Value = 9656.053
$ws.Range("I122").Value = 11104.0625
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 33312.1875
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -30862.1875
$ws.Range("N122").Value = -10700.0002
$ws.Range("H132").Value = 10110594
$ws.Range("I132").Value = 17552008
$ws.Range("K132").Value = 52656024
$ws.Range("M132").Value = -52653494
$ws.Range("H137").Value = 1213.2821
$ws.Range("I137").Value = 817.6818
$ws.Range("J137").Value = 1725.2354
$ws.Range("K137").Value = 2453.0454
$ws.Range("L137").Value = 5175.706200000001
$ws.Range("M137").Value = 96.95460000000003
$ws.Range("N137").Value = -10275.7062
$ws.Range("H138").Value = 440972.34
$ws.Range("J138").Value = 519373.97
$ws.Range("L138").Value = 1558121.91
$ws.Range("N138").Value = -1568401.91

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4551.552
$ws.Range("I32").Value = 3990.585
$ws.Range("K32").Value = 3990.585
$ws.Range("M32").Value = -3703.585
$ws.Range("H37").Value = 19919
$ws.Range("J37").Value = 20038
$ws.Range("L37").Value = 20038
$ws.Range("N37").Value = -20584
$ws.Range("H132").Value = 1158.3539
$ws.Range("I132").Value = 868.8868
$ws.Range("K132").Value = 2606.6604
$ws.Range("M132").Value = -76.66039999999975

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25806.334
$ws.Range("I82").Value = 16752.334
$ws.Range("J82").Value = 30333.334
$ws.Range("K82").Value = 16752.334
$ws.Range("L82").Value = 30333.334
$ws.Range("M82").Value = -16369.334
$ws.Range("N82").Value = -31099.334
$ws.Range("H85").Value = 25806.334
$ws.Range("I85").Value = 16752.334
$ws.Range("J85").Value = 30333.334
$ws.Range("K85").Value = 16752.334
$ws.Range("L85").Value = 30333.334
$ws.Range("M85").Value = -15426.334
$ws.Range("N85").Value = -32985.334
$ws.Range("H86").Value = 6163.375
$ws.Range("I86").Value = 6163.375
$ws.Range("K86").Value = 6163.375
$ws.Range("M86").Value = -5040.375
$ws.Range("H89").Value = 6163.375
$ws.Range("I89").Value = 6163.375
$ws.Range("K89").Value = 30816.875
$ws.Range("M89").Value = -25200.875
$ws.Range("H105").Value = 111113350
$ws.Range("I105").Value = 142859460
$ws.Range("K105").Value = 142859460
$ws.Range("M105").Value = -142857713

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1446.1111
$ws.Range("I31").Value = 1446.1111
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1446.1111
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1151.1111
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1446.1111
$ws.Range("I34").Value = 1446.1111
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1446.1111
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1244.1111
$ws.Range("N34").ClearContents()
$ws.Range("H134").Value = 649.9737
$ws.Range("I134").Value = 582.82855
$ws.Range("J134").Value = 1433.3334
$ws.Range("K134").Value = 1748.48565
$ws.Range("L134").Value = 4300.0002
$ws.Range("M134").Value = 786.5143500000001
$ws.Range("N134").Value = -9370.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1299.0303
$ws.Range("I5").Value = 1579.2727
$ws.Range("J5").Value = 738.5455
$ws.Range("K5").Value = 4737.8181
$ws.Range("L5").Value = 2215.6365
$ws.Range("M5").Value = -4625.8181
$ws.Range("N5").Value = -2439.6365
$ws.Range("H56").Value = 6058.9614
$ws.Range("I56").Value = 6058.9614
$ws.Range("K56").Value = 6058.9614
$ws.Range("M56").Value = -5528.9614
$ws.Range("H113").Value = 681.91895
$ws.Range("J113").Value = 684.19446
$ws.Range("L113").Value = 2052.58338
$ws.Range("N113").Value = -6392.58338
$ws.Range("H116").Value = 3500
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 10500
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -17384
$ws.Range("H117").Value = 1079.5
$ws.Range("I117").Value = 625.6667
$ws.Range("J117").Value = 1533.3334
$ws.Range("K117").Value = 1877.0001
$ws.Range("L117").Value = 4600.0002
$ws.Range("M117").Value = 1564.9999
$ws.Range("N117").Value = -11484.0002
$ws.Range("H118").Value = 800
$ws.Range("I118").Value = 800
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2400
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -1157
$ws.Range("N118").ClearContents()
$ws.Range("H131").Value = 21277798
$ws.Range("J131").Value = 1498.4166
$ws.Range("L131").Value = 4495.2498
$ws.Range("N131").Value = -14575.2498
$ws.Range("H135").Value = 1299.0303
$ws.Range("I135").Value = 1579.2727
$ws.Range("J135").Value = 738.5455
$ws.Range("K135").Value = 14213.4543
$ws.Range("L135").Value = 6646.9095
$ws.Range("M135").Value = -11678.4543
$ws.Range("N135").Value = -11716.9095

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 221.66667
$ws.Range("I2").Value = 237.25
$ws.Range("K2").Value = 237.25
$ws.Range("M2").Value = -124.25
$ws.Range("H28").Value = 5000
$ws.Range("J28").Value = 5000
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5384
$ws.Range("H126").Value = 2100
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -2930

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5344
$ws.Range("H55").Value = 345.45456
$ws.Range("I55").Value = 253.66667
$ws.Range("J55").Value = 542.1429
$ws.Range("K55").Value = 253.66667
$ws.Range("L55").Value = 542.1429
$ws.Range("M55").Value = -80.66667000000001
$ws.Range("N55").Value = -888.1429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33339254
$ws.Range("I62").Value = 35720308
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 35720308
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -35719684
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 33339254
$ws.Range("I65").Value = 35720308
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 178601540
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -178598420
$ws.Range("N65").Value = -28740
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H100").Value = 890.7692
$ws.Range("I100").Value = 1067.625
$ws.Range("J100").Value = 607.8
$ws.Range("K100").Value = 2135.25
$ws.Range("L100").Value = 1215.6
$ws.Range("M100").Value = -1594.25
$ws.Range("N100").Value = -2297.6
$ws.Range("H107").Value = 516.73334
$ws.Range("I107").Value = 495.91666
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1487.74998
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 432.2500199999999
$ws.Range("N107").Value = -5640
$ws.Range("H136").Value = 458.69232
$ws.Range("I136").Value = 220.61111
$ws.Range("K136").Value = 661.8333299999999
$ws.Range("M136").Value = 1888.16667

Write-Host "Applied 220 cell updates across 8 sheets"
